$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.05070303331912385
$ws.Range("D2").Value = 0.2332309276353044
$ws.Range("E2").Value = 0.04604448509153691
$ws.Range("F2").Value = 14.91324311680114
$ws.Range("G2").Value = 0.002918073986829643
$ws.Range("I2").Value = 11.19940568461692
$ws.Range("J2").Value = 0.4134903892799997
$ws.Range("L2").Value = 0.2071373558204996

# Row 3
$ws.Range("C3").Value = 0.04543810182379104
$ws.Range("D3").Value = 0.2199300229884784
$ws.Range("E3").Value = 0.04620450570436918
$ws.Range("F3").Value = 14.71638319715089
$ws.Range("G3").Value = 0.00293686367383339
$ws.Range("I3").Value = 11.04696120098765
$ws.Range("J3").Value = 0.4132330119168799
$ws.Range("L3").Value = 0.2084080281731815

# Row 4
$ws.Range("C4").Value = 0.04225366310134859
$ws.Range("D4").Value = 0.2119966071513772
$ws.Range("E4").Value = 0.04630899095280894
$ws.Range("F4").Value = 14.60697997074305
$ws.Range("G4").Value = 0.002948941771217428
$ws.Range("I4").Value = 10.96202227129328
$ws.Range("J4").Value = 0.4133321067958491
$ws.Range("L4").Value = 0.2092837934837455

# Row 5
$ws.Range("C5").Value = 0.04096734721755979
$ws.Range("D5").Value = 0.2088201796352394
$ws.Range("E5").Value = 0.04635313975012378
$ws.Range("F5").Value = 14.56523231834899
$ws.Range("G5").Value = 0.002954000779420692
$ws.Range("I5").Value = 10.92954969324299
$ws.Range("J5").Value = 0.4134365316262461
$ws.Range("L5").Value = 0.2096646463385952

# Row 6
$ws.Range("C6").Value = 0.04075442023054165
$ws.Range("D6").Value = 0.2082960865209884
$ws.Range("E6").Value = 0.0463605655552386
$ws.Range("F6").Value = 14.55846998389421
$ws.Range("G6").Value = 0.002954849132775558
$ws.Range("I6").Value = 10.92428587027433
$ws.Range("J6").Value = 0.4134577218993272
$ws.Range("L6").Value = 0.2097293328201211

# Row 7
$ws.Range("C7").Value = 0.04223627035560185
$ws.Range("D7").Value = 0.2119535427892458
$ws.Range("E7").Value = 0.04630957999683938
$ws.Range("F7").Value = 14.60640553169867
$ws.Range("G7").Value = 0.002949009442164172
$ws.Range("I7").Value = 10.96157571557259
$ws.Range("J7").Value = 0.4133332565880963
$ws.Range("L7").Value = 0.2092888328141953

# Row 8
$ws.Range("C8").Value = 0.04887720126285444
$ws.Range("D8").Value = 0.2285949924485067
$ws.Range("E8").Value = 0.04609836862103078
$ws.Range("F8").Value = 14.84295510381554
$ws.Range("G8").Value = 0.002924440967947825
$ws.Range("I8").Value = 11.14502178447378
$ws.Range("J8").Value = 0.4133479004541414
$ws.Range("L8").Value = 0.207555611077936

# Row 9
$ws.Range("C9").Value = 0.06231916247040203
$ws.Range("D9").Value = 0.263184562746801
$ws.Range("E9").Value = 0.04573349430180862
$ws.Range("F9").Value = 15.40017927948429
$ws.Range("G9").Value = 0.002880509321573744
$ws.Range("I9").Value = 11.57528781273408
$ws.Range("J9").Value = 0.4154463219016264
$ws.Range("L9").Value = 0.2049180654706717

# Row 10
$ws.Range("C10").Value = 0.07250352962647355
$ws.Range("D10").Value = 0.2899419851922005
$ws.Range("E10").Value = 0.04549529974900057
$ws.Range("F10").Value = 15.86993682970211
$ws.Range("G10").Value = 0.002850755126877473
$ws.Range("I10").Value = 11.93704484257205
$ws.Range("J10").Value = 0.4182935905190561
$ws.Range("L10").Value = 0.2034490359932306

# Row 11
$ws.Range("C11").Value = 0.07721590248584675
$ws.Range("D11").Value = 0.3024417553290277
$ws.Range("E11").Value = 0.04539339087049132
$ws.Range("F11").Value = 16.09755312645615
$ws.Range("G11").Value = 0.002837751999361088
$ws.Range("I11").Value = 12.11214286542548
$ws.Range("J11").Value = 0.4198826588740872
$ws.Range("L11").Value = 0.2028836782580399

# Row 12
$ws.Range("C12").Value = 0.07901282628145623
$ws.Range("D12").Value = 0.3072252256631884
$ws.Range("E12").Value = 0.04535572525306319
$ws.Range("F12").Value = 16.18581648607648
$ws.Range("G12").Value = 0.002832903357728166
$ws.Range("I12").Value = 12.18001544486094
$ws.Range("J12").Value = 0.4205275429966377
$ws.Range("L12").Value = 0.2026844950553865

# Row 13
$ws.Range("C13").Value = 0.07862525599406922
$ws.Range("D13").Value = 0.3061927437985332
$ws.Range("E13").Value = 0.04536379610536745
$ws.Range("F13").Value = 16.16671420453918
$ws.Range("G13").Value = 0.002833944267349885
$ws.Range("I13").Value = 12.16532732895189
$ws.Range("J13").Value = 0.4203867224229043
$ws.Range("L13").Value = 0.2027267280018847

# Row 14
$ws.Range("C14").Value = 0.07736348091361833
$ws.Range("D14").Value = 0.3028342731927012
$ws.Range("E14").Value = 0.04539027356872749
$ws.Range("F14").Value = 16.10477271695601
$ws.Range("G14").Value = 0.002837351595065305
$ws.Range("I14").Value = 12.11769507480551
$ws.Range("J14").Value = 0.4199348440661481
$ws.Range("L14").Value = 0.2028669919647896

# Row 15
$ws.Range("C15").Value = 0.07659225993901941
$ws.Range("D15").Value = 0.3007837249942327
$ws.Range("E15").Value = 0.04540661219376085
$ws.Range("F15").Value = 16.06710342733237
$ws.Range("G15").Value = 0.00283944846124351
$ws.Range("I15").Value = 12.08872457605668
$ws.Range("J15").Value = 0.4196637005738211
$ws.Range("L15").Value = 0.2029548521472933

# Row 16
$ws.Range("C16").Value = 0.0721972502902446
$ws.Range("D16").Value = 0.2891319392111598
$ws.Range("E16").Value = 0.04550208928828159
$ws.Range("F16").Value = 15.85534744512705
$ws.Range("G16").Value = 0.002851615522191988
$ws.Range("I16").Value = 11.92581811945212
$ws.Range("J16").Value = 0.4181957283883264
$ws.Range("L16").Value = 0.2034880632435829

# Row 17
$ws.Range("C17").Value = 0.06952211806694208
$ws.Range("D17").Value = 0.2820699329582794
$ws.Range("E17").Value = 0.04556231113997677
$ws.Range("F17").Value = 15.72905494878472
$ws.Range("G17").Value = 0.002859215109056289
$ws.Range("I17").Value = 11.82861426732711
$ws.Range("J17").Value = 0.4173710468894569
$ws.Range("L17").Value = 0.2038416025146788

# Row 18
$ws.Range("C18").Value = 0.0679908670333873
$ws.Range("D18").Value = 0.2780387114119094
$ws.Range("E18").Value = 0.0455975560892723
$ws.Range("F18").Value = 15.65772089075153
$ws.Range("G18").Value = 0.002863636355858848
$ws.Range("I18").Value = 11.77369331330897
$ws.Range("J18").Value = 0.4169243127738582
$ws.Range("L18").Value = 0.204054628540888

# Row 19
$ws.Range("C19").Value = 0.06747365471360922
$ws.Range("D19").Value = 0.2766789867769717
$ws.Range("E19").Value = 0.04560959372903195
$ws.Range("F19").Value = 15.63379069995233
$ws.Range("G19").Value = 0.002865141965306444
$ws.Range("I19").Value = 11.75526618979319
$ws.Range("J19").Value = 0.4167777705747255
$ws.Range("L19").Value = 0.2041284148469487

# Row 20
$ws.Range("C20").Value = 0.06980611584101837
$ws.Range("D20").Value = 0.2828185016634279
$ws.Range("E20").Value = 0.04555583762098481
$ws.Range("F20").Value = 15.74236334751271
$ws.Range("G20").Value = 0.002858400936718741
$ws.Range("I20").Value = 11.83885914654144
$ws.Range("J20").Value = 0.4174559728007523
$ws.Range("L20").Value = 0.2038029651063695

# Row 21
$ws.Range("C21").Value = 0.07773374851944936
$ws.Range("D21").Value = 0.3038193539927931
$ws.Range("E21").Value = 0.0453824714050266
$ws.Range("F21").Value = 16.12290969293542
$ws.Range("G21").Value = 0.002836348744664097
$ws.Range("I21").Value = 12.13164287197355
$ws.Range("J21").Value = 0.4200663936178302
$ws.Range("L21").Value = 0.2028253876129824

# Row 22
$ws.Range("C22").Value = 0.08298799772205712
$ws.Range("D22").Value = 0.3178380017773748
$ws.Range("E22").Value = 0.04527455774219735
$ws.Range("F22").Value = 16.383718559496
$ws.Range("G22").Value = 0.00282237505510566
$ws.Range("I22").Value = 12.3321524192009
$ws.Range("J22").Value = 0.4220243461843012
$ws.Range("L22").Value = 0.2022734206633103

# Row 23
$ws.Range("C23").Value = 0.08017666703742066
$ws.Range("D23").Value = 0.3103281254270769
$ws.Range("E23").Value = 0.04533166061554095
$ws.Range("F23").Value = 16.24338935905905
$ws.Range("G23").Value = 0.002829793331325822
$ws.Range("I23").Value = 12.22428075172303
$ws.Range("J23").Value = 0.4209559959259366
$ws.Range("L23").Value = 0.202560023546269

# Row 24
$ws.Range("C24").Value = 0.06967769959845782
$ws.Range("D24").Value = 0.2824799842057075
$ws.Range("E24").Value = 0.04555876235946599
$ws.Range("F24").Value = 15.73634265988682
$ws.Range("G24").Value = 0.002858768861461089
$ws.Range("I24").Value = 11.83422444154183
$ws.Range("J24").Value = 0.4174174925602756
$ws.Range("L24").Value = 0.2038204026504573

# Row 25
$ws.Range("C25").Value = 0.05863254563045928
$ws.Range("D25").Value = 0.2536024829852295
$ws.Range("E25").Value = 0.04582694284855382
$ws.Range("F25").Value = 15.23909464849731
$ws.Range("G25").Value = 0.002891946121025502
$ws.Range("I25").Value = 11.45107342679364
$ws.Range("J25").Value = 0.4154463219016264
$ws.Range("L25").Value = 0.2055496314700775
